$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Datos actualizados" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 28 de Mayo de 2020 a las 05:35"

# --- Update India row (row 13): new cases pushed totals/active up ---
$ws.Range("B13").Value = 158333
$ws.Range("C13").Value = 247
$ws.Range("E13").Value = 86050

# --- Update provincias Spain tail section (rows 160-162) ---
# Mongolia received +13 new cases, overtaking Guadalupe and Gibraltar in the
# descending sort by "Casos totales", so the three rows shift: Mongolia moves
# into row 160, Guadalupe moves down to row 161, Gibraltar moves down to row 162.
$ws.Range("A160").Value = "Mongolia"
$ws.Range("B160").Value = 161
$ws.Range("C160").Value = 13
$ws.Range("D160").Value = 43
$ws.Range("E160").Value = 118
$ws.Range("H160").Value = 0

$ws.Range("A161").Value = "Guadalupe"
$ws.Range("B161").Value = 161
$ws.Range("C161").Value = 0
$ws.Range("D161").Value = 115
$ws.Range("E161").Value = 32
$ws.Range("H161").Value = 14

$ws.Range("A162").Value = "Gibraltar"
$ws.Range("B162").Value = 157
$ws.Range("C162").Value = 0
$ws.Range("D162").Value = 147
$ws.Range("E162").Value = 10
$ws.Range("H162").Value = 0
